# EventSummary_Bags of Joy Distribution.xlsx - data corrections
# (commit: "Performance testing using NBench")
#
# 1. "Participated" sheet (sheet1): update Rating + free-text feedback for
#    the two existing respondent rows.
# 2. "NotParticipated" sheet (sheet2): no change.
# 3. "UnRegistered" sheet (sheet3): the single data row (row 2) is removed,
#    leaving just the header row.

$wb = $excel.ActiveWorkbook

$wsParticipated = $wb.Worksheets.Item("Participated")

# Row 2 (Associate ID 330721)
$wsParticipated.Range("F2").Value = 3
$wsParticipated.Range("G2").Value = "Nice Event"
$wsParticipated.Range("H2").Value = "Nothing Such"

# Row 3 (Associate ID 741602)
$wsParticipated.Range("F3").Value = 2
$wsParticipated.Range("G3").Value = "Okay"
$wsParticipated.Range("H3").Value = "Everything"

# "UnRegistered" sheet: drop the lone data row so only the header remains.
$wsUnRegistered = $wb.Worksheets.Item("UnRegistered")
$wsUnRegistered.Rows.Item(2).Delete()
